$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.474.45"
$ws.Range("E2").Value = "  -3.67%  "

$ws.Range("D3").Value = "1.991.28"
$ws.Range("E3").Value = "  -5.10%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'240.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.83%  "

$ws.Range("D6").Value = "'0.631"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.59%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "'56.38"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.66%  "

$ws.Range("D9").Value = "'59.06"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.21%  "

$ws.Range("E10").Value = "  -3.37%  "

$ws.Range("D11").Value = "'0.0724"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.41%  "

$ws.Range("E12").Value = "  -6.26%  "

$ws.Range("D13").Value = "'0.893"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.29%  "

$ws.Range("D14").Value = "'14.35"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.78%  "

$ws.Range("D15").Value = "2.278.81"
$ws.Range("E15").Value = "  -5.28%  "

$ws.Range("D16").Value = "'5.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.77%  "

$ws.Range("D17").Value = "1.990.07"
$ws.Range("E17").Value = "  -5.03%  "

$ws.Range("E18").Value = "  -1.10%  "

$ws.Range("D19").Value = "35.470.60"
$ws.Range("E19").Value = "  -3.67%  "

$ws.Range("D20").Value = "'69.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.12%  "

$ws.Range("D21").Value = "0.0₃0834"
$ws.Range("E21").Value = "  -5.53%  "

$ws.Range("D22").Value = "'231.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.34%  "

$ws.Range("E23").Value = "  -8.52%  "

$ws.Range("E24").Value = "  +0.14%  "

$ws.Range("D25").Value = "'2.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -6.03%  "

$ws.Range("E26").Value = "  +4.41%  "

$ws.Range("D27").Value = "'9.11"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.74%  "

$ws.Range("D28").Value = "'162.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.64%  "

$ws.Range("E29").Value = "  -7.60%  "

$ws.Range("E30").Value = "  -3.71%  "

$ws.Range("E31").Value = "  -2.46%  "

$ws.Range("D32").Value = "'4.76"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -9.35%  "

$ws.Range("D33").Value = "'0.0584"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.04%  "

$ws.Range("D34").Value = "'0.0902"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.76%  "

$ws.Range("E35").Value = "  -10.58%  "

$ws.Range("E36").Value = "  +0.05%  "

$ws.Range("E37").Value = "  -8.72%  "

$ws.Range("E38").Value = "  -2.48%  "

$ws.Range("E39").Value = "  -0.90%  "

$ws.Range("E40").Value = "  -7.71%  "

$ws.Range("E41").Value = "  -0.94%  "

$ws.Range("E42").Value = "  -5.64%  "

$ws.Range("E43").Value = "  -7.03%  "

$ws.Range("D44").Value = "'0.0883"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.20%  "

$ws.Range("D45").Value = "'90.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.86%  "

$ws.Range("D46").Value = "1.366.68"
$ws.Range("E46").Value = "  -3.70%  "

$ws.Range("E47").Value = "  -6.08%  "

$ws.Range("D48").Value = "'15.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.43%  "

$ws.Range("E49").Value = "  -0.71%  "

$ws.Range("D50").Value = "'2.25"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.90%  "

$ws.Range("D51").Value = "'45.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.52%  "
